$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.361.87"
$ws.Range("E2").Value = "  -1.47%  "
$ws.Range("D3").Value = "2.594.55"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.54"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.90"
$ws.Range("E6").Value = "  -1.58%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.72"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.384"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.49"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").Value = "3.063.11"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "63.227.00"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000155"
$ws.Range("E16").Value = "  +3.56%  "
$ws.Range("D17").Value = "2.600.73"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.02"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.67"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.52"
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.81"
$ws.Range("E21").Value = "  -2.73%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.54"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  -2.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.12"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.64"
$ws.Range("E26").Value = "  -3.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.22"
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "552.20"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  -3.49%  "
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("D32").Value = "0.0₃0847"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.75"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.25"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.47"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.410"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.34"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("E40").Value = "  -4.97%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "164.82"
$ws.Range("E42").Value = "  -4.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.98"
$ws.Range("E43").Value = "  +1.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.84"
$ws.Range("E44").Value = "  +5.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0576"
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.09"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0954"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.96"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").Value = "0.0₆0222"
$ws.Range("E51").Value = "  +11.28%  "
